$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.269999999999999
$ws.Range("A8").Value = -21.107
$ws.Range("A10").Value = -20.973
$ws.Range("A12").Value = -21.452
$ws.Range("B13").Value = 7.115
$ws.Range("A18").Value = -21.766
$ws.Range("D20").Value = -8.222
$ws.Range("A25").Value = -21.534
